$d = $word.ActiveDocument

# --- Insert the new "Introduction" section before the trailing bookmark paragraph ---
# The document's final paragraph is an (otherwise empty) paragraph that only carries the
# _GoBack bookmark. We insert the two new paragraphs immediately before it, so the
# bookmark paragraph remains the last "real" paragraph, as in the target document.
$bookmarkPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertPoint = $bookmarkPara.Range
$insertPoint.Collapse(1)  # wdCollapseStart

$insertPoint.InsertBefore("Introduction`r")
$introPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$introPara.Range.Style = $d.Styles.Item("Heading 2")

$bookmarkPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertPoint2 = $bookmarkPara.Range
$insertPoint2.Collapse(1)
$bodyText = "This report explains how to use GitHub for version control. It includes steps for creating repositories, adding collaborators, making commits, and managing branches. The purpose is to understand how team collaboration works on GitHub."
$insertPoint2.InsertBefore($bodyText + "`r")

# --- Append a trailing empty paragraph after the (still-last) bookmark paragraph ---
$endRange = $d.Content
$endRange.Collapse(0)  # wdCollapseEnd
$endRange.InsertParagraphAfter()

Write-Output ("Final paragraph count: " + $d.Paragraphs.Count)
